# Apply the edits described by the diff:
#  - C8 changes from "物品" to "數量"
#  - C13 changes from "數量" to "單位"
#  - Active selection moves from B11 to C8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = "數量"
$ws.Range("C13").Value = "單位"

$ws.Range("C8").Select()
